$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 was mistakenly stored as text "4" - correct it to a real number 4
$ws.Range("B8").Value = 4

# Append a new annotation row (row 9) for Ying Tang
$ws.Range("A9").Value = "Ying Tang"

# B9 holds the score "2" as text (matching the style of the original rows
# before they were later normalized to numbers), so force text formatting
# before assigning the value to avoid Excel auto-converting it to a number.
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2"
$ws.Range("B9").Style = "Normal"

$ws.Range("C9").Value = " I do not enjoy,convoluted"
$ws.Range("D9").Value = "CRT"
$ws.Range("E9").Value = "WRI"
$ws.Range("F9").Value = "f5b44bd7-9311-4cfc-b939-3b86c20706ac"
$ws.Range("G9").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H9").Value = "On top of this, I do not enjoy the style the paper is written in, the language is convoluted."
